$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N6").ClearContents()
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -2739.5
$ws.Range("H6").Value = 950.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2851.5
$ws.Range("I6").Value = 950.5
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("J75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("H76").Value = 3369.5652
$ws.Range("J76").Value = 3550
$ws.Range("M76").Value = -3016.5789
$ws.Range("K76").Value = 3331.5789
$ws.Range("L76").Value = 3550
$ws.Range("N76").Value = -4180
$ws.Range("I76").Value = 3331.5789
$ws.Range("N78").ClearContents()
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("K79").Value = 3331.5789
$ws.Range("L79").Value = 3550
$ws.Range("I79").Value = 3331.5789
$ws.Range("N79").Value = -5734
$ws.Range("M79").Value = -2239.5789
$ws.Range("H79").Value = 3369.5652
$ws.Range("J79").Value = 3550
$ws.Range("I80").Value = 315.66666
$ws.Range("H80").Value = 467.26923
$ws.Range("K80").Value = 946.9999799999999
$ws.Range("L80").Value = 3312
$ws.Range("J80").Value = 1104
$ws.Range("N80").Value = -5308
$ws.Range("M80").Value = 51.00002000000006
$ws.Range("L83").Value = 9936
$ws.Range("N83").Value = -19920
$ws.Range("J83").Value = 1104
$ws.Range("M83").Value = 2151.00006
$ws.Range("K83").Value = 2840.99994
$ws.Range("I83").Value = 315.66666
$ws.Range("H83").Value = 467.26923
$ws.Range("J87").Value = 20787.8
$ws.Range("N87").Value = -23283.8
$ws.Range("L87").Value = 20787.8
$ws.Range("H87").Value = 20787.8
$ws.Range("H90").Value = 20787.8
$ws.Range("L90").Value = 62363.39999999999
$ws.Range("J90").Value = 20787.8
$ws.Range("N90").Value = -74843.39999999999
$ws.Range("M94").Value = -3121.2727
$ws.Range("H94").Value = 3572.2727
$ws.Range("I94").Value = 3572.2727
$ws.Range("K94").Value = 3572.2727
$ws.Range("J112").Value = 1365
$ws.Range("L112").Value = 4095
$ws.Range("N112").Value = -6311
$ws.Range("H112").Value = 1351.5476
$ws.Range("K113").Value = 1747.5
$ws.Range("N113").Value = -15119.333
$ws.Range("H113").Value = 7363.364
$ws.Range("M113").Value = 1506.5
$ws.Range("I113").Value = 1747.5
$ws.Range("J113").Value = 8611.333000000001
$ws.Range("L113").Value = 8611.333000000001
$ws.Range("L120").Value = 30000
$ws.Range("H120").Value = 30000
$ws.Range("N120").Value = -39676
$ws.Range("J120").Value = 30000
$ws.Range("H123").Value = 40910
$ws.Range("J123").Value = 40910
$ws.Range("N123").Value = -50710
$ws.Range("L123").Value = 40910
$ws.Range("J132").Value = 1114097.8
$ws.Range("M132").Value = -150606260
$ws.Range("N132").Value = -3347353.4
$ws.Range("H132").Value = 34968464
$ws.Range("K132").Value = 150608790
$ws.Range("L132").Value = 3342293.4
$ws.Range("I132").Value = 50202930
$ws.Range("I138").Value = 1666.4706
$ws.Range("N138").Value = -20272.7555
$ws.Range("K138").Value = 4999.4118
$ws.Range("H138").Value = 2902.197
$ws.Range("M138").Value = 140.5882000000001
$ws.Range("L138").Value = 9992.755500000001
$ws.Range("J138").Value = 3330.9185

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L32").Value = 5429.1665
$ws.Range("N32").Value = -6003.1665
$ws.Range("M32").Value = -4431.4443
$ws.Range("I32").Value = 4718.4443
$ws.Range("H32").Value = 4868.0703
$ws.Range("J32").Value = 5429.1665
$ws.Range("K32").Value = 4718.4443
$ws.Range("I45").Value = 4337.3335
$ws.Range("M45").Value = -3960.3335
$ws.Range("K45").Value = 4337.3335
$ws.Range("H45").Value = 2737.8333
$ws.Range("M74").Value = -380326.06
$ws.Range("I74").Value = 381200.06
$ws.Range("N74").Value = -3589
$ws.Range("H74").Value = 237886.64
$ws.Range("J74").Value = 1841
$ws.Range("L74").Value = 1841
$ws.Range("K74").Value = 381200.06
$ws.Range("J77").Value = 1841
$ws.Range("L77").Value = 9205
$ws.Range("K77").Value = 1906000.3
$ws.Range("M77").Value = -1901632.3
$ws.Range("H77").Value = 237886.64
$ws.Range("I77").Value = 381200.06
$ws.Range("N77").Value = -17941
$ws.Range("J132").Value = 4631.1
$ws.Range("M132").Value = -2735.5712
$ws.Range("N132").Value = -18953.3
$ws.Range("H132").Value = 2682.9033
$ws.Range("K132").Value = 5265.5712
$ws.Range("L132").Value = 13893.3
$ws.Range("I132").Value = 1755.1904
$ws.Range("N137").Value = -50511
$ws.Range("L137").Value = 40311
$ws.Range("J137").Value = 40311
$ws.Range("H137").Value = 40311

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J134").Value = 5693.1787
$ws.Range("N134").Value = -22149.5361
$ws.Range("I134").Value = 1328.6957
$ws.Range("L134").Value = 17079.5361
$ws.Range("H134").Value = 3724.8823
$ws.Range("K134").Value = 3986.0871
$ws.Range("M134").Value = -1451.0871

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 676610.7
$ws.Range("N31").Value = -3916.4849
$ws.Range("M31").Value = -676315.7
$ws.Range("L31").Value = 3326.4849
$ws.Range("H31").Value = 257396
$ws.Range("K31").Value = 676610.7
$ws.Range("J31").Value = 3326.4849
$ws.Range("L34").Value = 3326.4849
$ws.Range("J34").Value = 3326.4849
$ws.Range("N34").Value = -3730.4849
$ws.Range("M34").Value = -676408.7
$ws.Range("K34").Value = 676610.7
$ws.Range("I34").Value = 676610.7
$ws.Range("H34").Value = 257396
$ws.Range("J106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("L106").Value = 0
$ws.Range("H106").Value = 0
$ws.Range("I115").Value = 24899
$ws.Range("H115").Value = 24899.5
$ws.Range("K115").Value = 24899
$ws.Range("M115").Value = -23724

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N68").Value = -11827.76
$ws.Range("I68").Value = 774.3333
$ws.Range("J68").Value = 3401.92
$ws.Range("H68").Value = 2416.575
$ws.Range("M68").Value = -1511.9999
$ws.Range("L68").Value = 10205.76
$ws.Range("K68").Value = 2322.9999
$ws.Range("J71").Value = 3401.92
$ws.Range("M71").Value = -2912.9997
$ws.Range("L71").Value = 30617.28
$ws.Range("N71").Value = -38729.28
$ws.Range("H71").Value = 2416.575
$ws.Range("K71").Value = 6968.9997
$ws.Range("I71").Value = 774.3333
$ws.Range("J107").Value = 24113.408
$ws.Range("N107").Value = -76180.224
$ws.Range("I107").Value = 429.14706
$ws.Range("K107").Value = 1287.44118
$ws.Range("L107").Value = 72340.224
$ws.Range("H107").Value = 13789.5
$ws.Range("M107").Value = 632.55882
$ws.Range("K113").Value = 1992.5001
$ws.Range("N113").Value = -37507541
$ws.Range("H113").Value = 4465093.5
$ws.Range("M113").Value = 177.4999
$ws.Range("I113").Value = 664.1667
$ws.Range("J113").Value = 12501067
$ws.Range("L113").Value = 37503201
$ws.Range("L122").Value = 27046.8936
$ws.Range("J122").Value = 3005.2104
$ws.Range("N122").Value = -31946.8936
$ws.Range("H122").Value = 2247.1228
$ws.Range("J131").Value = 855.59375
$ws.Range("M131").Value = 3379.5
$ws.Range("H131").Value = 843.51
$ws.Range("K131").Value = 1660.5
$ws.Range("N131").Value = -12646.78125
$ws.Range("L131").Value = 2566.78125
$ws.Range("I131").Value = 553.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L122").Value = 24990
$ws.Range("I122").Value = 9999
$ws.Range("J122").Value = 8330
$ws.Range("M122").Value = -27547
$ws.Range("K122").Value = 29997
$ws.Range("N122").Value = -29890
$ws.Range("H122").Value = 8747.25
$ws.Range("J132").Value = 4134.6
$ws.Range("M132").Value = -2765
$ws.Range("N132").Value = -17463.8
$ws.Range("H132").Value = 3293.7742
$ws.Range("K132").Value = 5295
$ws.Range("L132").Value = 12403.8
$ws.Range("I132").Value = 1765

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N18").Value = -20344
$ws.Range("L18").Value = 20000
$ws.Range("J18").Value = 20000
$ws.Range("H18").Value = 20000
$ws.Range("I46").Value = 653.4545000000001
$ws.Range("H46").Value = 1012.129
$ws.Range("K46").Value = 653.4545000000001
$ws.Range("M46").Value = -465.4545000000001
$ws.Range("J132").Value = 4919.0713
$ws.Range("M132").Value = -5662.6844
$ws.Range("N132").Value = -19817.2139
$ws.Range("H132").Value = 3659.2122
$ws.Range("K132").Value = 8192.6844
$ws.Range("L132").Value = 14757.2139
$ws.Range("I132").Value = 2730.8948
$ws.Range("M136").Value = -2365.928400000001
$ws.Range("N136").Value = -22658.8242
$ws.Range("L136").Value = 17558.8242
$ws.Range("J136").Value = 5852.9414
$ws.Range("H136").Value = 3949.7097
$ws.Range("K136").Value = 4915.928400000001
$ws.Range("I136").Value = 1638.6428

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I80").Value = 0
$ws.Range("H80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("J83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("K83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J107").Value = 978
$ws.Range("N107").Value = -6774
$ws.Range("I107").Value = 582.25
$ws.Range("K107").Value = 1746.75
$ws.Range("L107").Value = 2934
$ws.Range("H107").Value = 802.1111
$ws.Range("M107").Value = 173.25
$ws.Range("J132").Value = 3581.0908
$ws.Range("M132").Value = -1304.9231
$ws.Range("N132").Value = -15803.2724
$ws.Range("H132").Value = 2333.75
$ws.Range("K132").Value = 3834.9231
$ws.Range("L132").Value = 10743.2724
$ws.Range("I132").Value = 1278.3077
$ws.Range("M136").Value = -5380.7145
$ws.Range("N136").Value = -24047.0772
$ws.Range("L136").Value = 18947.0772
$ws.Range("J136").Value = 6315.6924
$ws.Range("H136").Value = 5030.45
$ws.Range("K136").Value = 7930.7145
$ws.Range("I136").Value = 2643.5715
